# Diagrama de Pareto.pptx edit script
# - Move the "Usos" slide (index 3) to index 5
# - Repurpose that slide (now at index 5) with new intro content
# - Insert two new slides ("Características principales", "¿Cuándo se utiliza?")
# - Update the footer date placeholder text across layouts + master
# - Add presentation slide guides (best effort)

$p = $ppt.ActivePresentation

# --- 1. Reorder: move the "Usos" slide from position 3 to position 5 ---------------
$usos = $p.Slides.Item(3)
$usos.MoveTo(5)

# --- 2. Repurpose the former "Usos" slide (now at position 5) ----------------------
$s5 = $p.Slides.Item(5)

$s5Title = $s5.Shapes.Item(1)
$s5Title.Name = "Título 1"
$s5Title.TextFrame.TextRange.Text = ""

$s5Body = $s5.Shapes.Item(2)
$s5Body.Name = "Marcador de contenido 2"
$s5Body.TextFrame.TextRange.Text = "Las Tablas y Diagramas de Pareto son herramientas de representación utilizadas para visualizar el Análisis de Pareto.`rEl Diagrama de Pareto es la representación gráfica de la Tabla de Pareto correspondiente.`r"

# --- 3. Insert new slide 6: "Características principales" -------------------------
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$s6 = $p.Slides.AddSlide(6, $layout)
$s6.Shapes.Item(1).Name = "Título 1"
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Características principales "

$s6.Shapes.Item(2).Name = "Marcador de contenido 2"
$s6Body = $s6.Shapes.Item(2).TextFrame.TextRange
$s6Body.Text = "Simplicidad: Tanto la Tabla como el Diagrama de Pareto no requieren ni cálculos complejos ni técnicas sofisticadas de representación gráfica. `rImpacto visual: El Diagrama de Pareto comunica de forma clara, evidente y de un ""vistazo"", el resultado del análisis de comparación y priorización."

# --- 4. Insert new slide 7: "¿Cuándo se utiliza?" -----------------------------------
$s7 = $p.Slides.AddSlide(7, $layout)
$s7.Shapes.Item(1).Name = "Título 1"
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "¿Cuándo se utiliza?"

$s7.Shapes.Item(2).Name = "Marcador de contenido 2"
$s7Lines = @(
    "Al identificar un producto o servicio para el análisis para mejorar la calidad.",
    "Cuando existe la necesidad de llamar la atención a los problema o causas de una forma sistemática.",
    "Al identificar oportunidades para mejorar",
    "Al analizar las diferentes agrupaciones de datos (ej: por producto, por segmento, del mercado, área geográfica, etc.)",
    "Al buscar las causas principales de los problemas y establecer la prioridad de las soluciones",
    "Al evaluar los resultados de los cambos efectuados a un proceso (antes y después)",
    "Cuando los datos puedan clasificarse en categorías",
    "Cuando el rango de cada categoría es importante"
)
$s7.Shapes.Item(2).TextFrame.TextRange.Text = [string]::Join("`r", $s7Lines)

# --- 5. Footer date placeholder: 17/04/2015 -> 21/04/2015 --------------------------
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $lyt = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $lyt.Shapes.Count; $j++) {
        $shp = $lyt.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "17/04/2015") {
                $shp.TextFrame.TextRange.Text = "21/04/2015"
            }
        }
    }
}
for ($j = 1; $j -le $p.SlideMaster.Shapes.Count; $j++) {
    $shp = $p.SlideMaster.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "17/04/2015") {
            $shp.TextFrame.TextRange.Text = "21/04/2015"
        }
    }
}

# --- 6. Slide guides (best effort; cosmetic only) -----------------------------------
try {
    $g1 = $p.Guides.Add(1, 2160)
    $g2 = $p.Guides.Add(2, 2880)
} catch {
}
